$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 3 (pushes all existing match rows down by one).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new match (Australia A-League: Melbourne City vs Melbourne Victory).
$ws.Cells.Item(3,1).Value = "U32gW6PH"  # A3
$ws.Cells.Item(3,2).Value = "26/10/2024"  # B3
$ws.Cells.Item(3,3).Value = "05:35"  # C3
$ws.Cells.Item(3,4).Value = "AUSTRALIA - A-LEAGUE"  # D3
$ws.Cells.Item(3,5).Value = "Melbourne City"  # E3
$ws.Cells.Item(3,6).Value = "Melbourne Victory"  # F3
$ws.Cells.Item(3,7).Value = 2.75  # G3
$ws.Cells.Item(3,8).Value = 3.5  # H3
$ws.Cells.Item(3,9).Value = 2.45  # I3
$ws.Cells.Item(3,10).Value = 3.25  # J3
$ws.Cells.Item(3,11).Value = 2.25  # K3
$ws.Cells.Item(3,12).Value = 3  # L3
$ws.Cells.Item(3,13).Value = 1.04  # M3
$ws.Cells.Item(3,14).Value = 13  # N3
$ws.Cells.Item(3,15).Value = 1.22  # O3
$ws.Cells.Item(3,16).Value = 4.33  # P3
$ws.Cells.Item(3,17).Value = 1.73  # Q3
$ws.Cells.Item(3,18).Value = 2.1  # R3
$ws.Cells.Item(3,19).Value = 1.33  # S3
$ws.Cells.Item(3,20).Value = 3.25  # T3
$ws.Cells.Item(3,21).Value = 1.62  # U3
$ws.Cells.Item(3,22).Value = 2.2  # V3
$ws.Cells.Item(3,23).Value = 11  # W3
$ws.Cells.Item(3,24).Value = 15  # X3
$ws.Cells.Item(3,25).Value = 10  # Y3
$ws.Cells.Item(3,26).Value = 29  # Z3
$ws.Cells.Item(3,27).Value = 21  # AA3
$ws.Cells.Item(3,28).Value = 26  # AB3
$ws.Cells.Item(3,29).Value = 13  # AC3
$ws.Cells.Item(3,30).Value = 7  # AD3
$ws.Cells.Item(3,31).Value = 12  # AE3
$ws.Cells.Item(3,32).Value = 41  # AF3
$ws.Cells.Item(3,33).Value = 151  # AG3
$ws.Cells.Item(3,34).Value = 10  # AH3
$ws.Cells.Item(3,35).Value = 13  # AI3
$ws.Cells.Item(3,36).Value = 9.5  # AJ3
$ws.Cells.Item(3,37).Value = 23  # AK3
$ws.Cells.Item(3,38).Value = 19  # AL3
$ws.Cells.Item(3,39).Value = 23  # AM3
$ws.Cells.Item(3,40).Value = 5  # AN3
$ws.Cells.Item(3,41).Value = 15  # AO3
$ws.Cells.Item(3,42).Value = 21  # AP3
$ws.Cells.Item(3,43).Value = 41  # AQ3
$ws.Cells.Item(3,44).Value = 51  # AR3
$ws.Cells.Item(3,45).Value = 126  # AS3
$ws.Cells.Item(3,46).Value = 3.25  # AT3
$ws.Cells.Item(3,47).Value = 7.5  # AU3
$ws.Cells.Item(3,48).Value = 41  # AV3
$ws.Cells.Item(3,49).Value = 451  # AW3
$ws.Cells.Item(3,50).Value = 4.75  # AX3
$ws.Cells.Item(3,51).Value = 13  # AY3
$ws.Cells.Item(3,52).Value = 21  # AZ3
$ws.Cells.Item(3,53).Value = 41  # BA3
$ws.Cells.Item(3,54).Value = 51  # BB3
$ws.Cells.Item(3,55).Value = 126  # BC3
$ws.Cells.Item(3,56).Value = 126  # BD3

# The match that is now on the last row (row 8, previously row 7 "Suwon FC - Seoul") was
# re-scraped with refreshed odds; update its values accordingly.
$ws.Cells.Item(8,1).Value = "86Td3Gio"  # A8
$ws.Cells.Item(8,2).Value = "26/10/2024"  # B8
$ws.Cells.Item(8,3).Value = "04:30"  # C8
$ws.Cells.Item(8,4).Value = "SOUTH KOREA - K LEAGUE 1"  # D8
$ws.Cells.Item(8,5).Value = "Suwon FC"  # E8
$ws.Cells.Item(8,6).Value = "Seoul"  # F8
$ws.Cells.Item(8,7).Value = 3.4  # G8
$ws.Cells.Item(8,8).Value = 3.3  # H8
$ws.Cells.Item(8,9).Value = 2.1  # I8
$ws.Cells.Item(8,10).Value = 3.75  # J8
$ws.Cells.Item(8,11).Value = 2.2  # K8
$ws.Cells.Item(8,12).Value = 2.75  # L8
$ws.Cells.Item(8,13).Value = 1.05  # M8
$ws.Cells.Item(8,14).Value = 11  # N8
$ws.Cells.Item(8,15).Value = 1.29  # O8
$ws.Cells.Item(8,16).Value = 3.5  # P8
$ws.Cells.Item(8,17).Value = 1.98  # Q8
$ws.Cells.Item(8,18).Value = 1.88  # R8
$ws.Cells.Item(8,19).Value = 1.4  # S8
$ws.Cells.Item(8,20).Value = 2.75  # T8
$ws.Cells.Item(8,21).Value = 1.73  # U8
$ws.Cells.Item(8,22).Value = 2  # V8
$ws.Cells.Item(8,23).Value = 11  # W8
$ws.Cells.Item(8,24).Value = 17  # X8
$ws.Cells.Item(8,25).Value = 12  # Y8
$ws.Cells.Item(8,26).Value = 34  # Z8
$ws.Cells.Item(8,27).Value = 26  # AA8
$ws.Cells.Item(8,28).Value = 34  # AB8
$ws.Cells.Item(8,29).Value = 10  # AC8
$ws.Cells.Item(8,30).Value = 6.5  # AD8
$ws.Cells.Item(8,31).Value = 13  # AE8
$ws.Cells.Item(8,32).Value = 41  # AF8
$ws.Cells.Item(8,33).Value = 201  # AG8
$ws.Cells.Item(8,34).Value = 8  # AH8
$ws.Cells.Item(8,35).Value = 10  # AI8
$ws.Cells.Item(8,36).Value = 9  # AJ8
$ws.Cells.Item(8,37).Value = 19  # AK8
$ws.Cells.Item(8,38).Value = 17  # AL8
$ws.Cells.Item(8,39).Value = 26  # AM8
$ws.Cells.Item(8,40).Value = 5.5  # AN8
$ws.Cells.Item(8,41).Value = 19  # AO8
$ws.Cells.Item(8,42).Value = 26  # AP8
$ws.Cells.Item(8,43).Value = 51  # AQ8
$ws.Cells.Item(8,44).Value = 81  # AR8
$ws.Cells.Item(8,45).Value = 151  # AS8
$ws.Cells.Item(8,46).Value = 2.75  # AT8
$ws.Cells.Item(8,47).Value = 7.5  # AU8
$ws.Cells.Item(8,48).Value = 51  # AV8
$ws.Cells.Item(8,49).Value = 501  # AW8
$ws.Cells.Item(8,50).Value = 4.33  # AX8
$ws.Cells.Item(8,51).Value = 12  # AY8
$ws.Cells.Item(8,52).Value = 21  # AZ8
$ws.Cells.Item(8,53).Value = 41  # BA8
$ws.Cells.Item(8,54).Value = 51  # BB8
$ws.Cells.Item(8,55).Value = 151  # BC8
$ws.Cells.Item(8,56).Value = 51  # BD8
